$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the parish/group id numbers in column B (rows 2-65) -----------------
# The ids were bogus (4000000000, a "supergroup" placeholder); replace them
# with the correct id (2000000000) for every data row.
for ($r = 2; $r -le 65; $r++) {
    $ws.Cells.Item($r, 2).Value = 2000000000
}

# --- Remove the two bogus / made-up rows at the bottom (66 & 67) -------------
# Their ids weren't possible (100000000558 / 0), and the corresponding shared
# strings ("ACTIVITES NON CLASSEES" / "__Trv ::Des1") are dropped too since
# nothing references them anymore.
$ws.Range("B66:D67").ClearContents() | Out-Null

# --- Restore the view: scroll down to where the edits were made, and leave
#     the selection on the last touched cell -----------------------------------
$excel.ActiveWindow.ScrollRow = 31
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E59").Select() | Out-Null
